$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 ("repaymentstrategy") value changes from "RBI (India)" to the new
# scenario string "Overdue/Due Fee/Int,Principal".
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Update the active selection to match the edited cell.
$ws.Range("B17").Select()
